$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2463200.8
$ws.Range("I6").Value = 2463200.8
$ws.Range("K6").Value = 7389602.399999999
$ws.Range("M6").Value = -7389490.399999999
$ws.Range("H38").Value = 320.75
$ws.Range("I38").Value = 327.66666
$ws.Range("J38").Value = 300
$ws.Range("K38").Value = 982.9999799999999
$ws.Range("L38").Value = 900
$ws.Range("M38").Value = -610.9999799999999
$ws.Range("N38").Value = -1644
$ws.Range("H39").Value = 2285.4614
$ws.Range("I39").Value = 1419.375
$ws.Range("J39").Value = 3671.2
$ws.Range("K39").Value = 4258.125
$ws.Range("L39").Value = 11013.6
$ws.Range("M39").Value = -3962.125
$ws.Range("N39").Value = -11605.6
$ws.Range("H70").Value = 77383450
$ws.Range("I70").Value = 41669100
$ws.Range("J70").Value = 104169200
$ws.Range("K70").Value = 125007300
$ws.Range("L70").Value = 312507600
$ws.Range("M70").Value = -125007030
$ws.Range("N70").Value = -312508140
$ws.Range("H73").Value = 77383450
$ws.Range("I73").Value = 41669100
$ws.Range("J73").Value = 104169200
$ws.Range("K73").Value = 125007300
$ws.Range("L73").Value = 312507600
$ws.Range("M73").Value = -125006364
$ws.Range("N73").Value = -312509472
$ws.Range("H80").Value = 7844663.5
$ws.Range("I80").Value = 13900211
$ws.Range("J80").Value = 58959.285
$ws.Range("K80").Value = 41700633
$ws.Range("L80").Value = 176877.855
$ws.Range("M80").Value = -41699635
$ws.Range("N80").Value = -178873.855
$ws.Range("H83").Value = 7844663.5
$ws.Range("I83").Value = 13900211
$ws.Range("J83").Value = 58959.285
$ws.Range("K83").Value = 125101899
$ws.Range("L83").Value = 530633.5650000001
$ws.Range("M83").Value = -125096907
$ws.Range("N83").Value = -540617.5650000001
$ws.Range("H98").Value = 35717630
$ws.Range("I98").Value = 40003236
$ws.Range("K98").Value = 40003236
$ws.Range("M98").Value = -40001738
$ws.Range("H122").Value = 35717630
$ws.Range("I122").Value = 40003236
$ws.Range("K122").Value = 120009708
$ws.Range("M122").Value = -120007258
$ws.Range("H132").Value = 2336.2856
$ws.Range("I132").Value = 2261.389
$ws.Range("J132").Value = 2471.1
$ws.Range("K132").Value = 6784.167
$ws.Range("L132").Value = 7413.299999999999
$ws.Range("M132").Value = -4254.167
$ws.Range("N132").Value = -12473.3
$ws.Range("H138").Value = 4185.615
$ws.Range("I138").Value = 1015
$ws.Range("J138").Value = 6903.2856
$ws.Range("K138").Value = 3045
$ws.Range("L138").Value = 20709.8568
$ws.Range("M138").Value = 2095
$ws.Range("N138").Value = -30989.8568
$ws.Range("H141").Value = 2325
$ws.Range("I141").Value = 2361.6667
$ws.Range("K141").Value = 7085.000100000001
$ws.Range("M141").Value = -1905.000100000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 65
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 70
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 42
$ws.Range("N5").Value = -274
$ws.Range("H32").Value = 1840364.9
$ws.Range("I32").Value = 1924990.1
$ws.Range("K32").Value = 1924990.1
$ws.Range("M32").Value = -1924703.1
$ws.Range("H61").Value = 3679.3242
$ws.Range("I61").Value = 2307.2727
$ws.Range("J61").Value = 14998.75
$ws.Range("K61").Value = 2307.2727
$ws.Range("L61").Value = 14998.75
$ws.Range("M61").Value = -2095.2727
$ws.Range("N61").Value = -15422.75
$ws.Range("H74").Value = 55035.676
$ws.Range("J74").Value = 4253.8184
$ws.Range("L74").Value = 4253.8184
$ws.Range("N74").Value = -6001.8184
$ws.Range("H77").Value = 55035.676
$ws.Range("J77").Value = 4253.8184
$ws.Range("L77").Value = 21269.092
$ws.Range("N77").Value = -30005.092
$ws.Range("H88").Value = 2443.238
$ws.Range("I88").Value = 1736.625
$ws.Range("J88").Value = 2878.077
$ws.Range("K88").Value = 1736.625
$ws.Range("L88").Value = 2878.077
$ws.Range("M88").Value = -1330.625
$ws.Range("N88").Value = -3690.077
$ws.Range("H91").Value = 2443.238
$ws.Range("I91").Value = 1736.625
$ws.Range("J91").Value = 2878.077
$ws.Range("K91").Value = 1736.625
$ws.Range("L91").Value = 2878.077
$ws.Range("M91").Value = -332.625
$ws.Range("N91").Value = -5686.077
$ws.Range("H110").Value = 1171.8334
$ws.Range("I110").Value = 1171.8334
$ws.Range("K110").Value = 1171.8334
$ws.Range("M110").Value = 873.1666
$ws.Range("H132").Value = 7764.2
$ws.Range("I132").Value = 6847.263
$ws.Range("J132").Value = 8853.0625
$ws.Range("K132").Value = 20541.789
$ws.Range("L132").Value = 26559.1875
$ws.Range("M132").Value = -18011.789
$ws.Range("N132").Value = -31619.1875
$ws.Range("H136").Value = 3679.3242
$ws.Range("I136").Value = 2307.2727
$ws.Range("J136").Value = 14998.75
$ws.Range("K136").Value = 6921.8181
$ws.Range("L136").Value = 44996.25
$ws.Range("M136").Value = -4371.8181
$ws.Range("N136").Value = -50096.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 65
$ws.Range("I4").Value = 70
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 70
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = -280
$ws.Range("H55").Value = 58688
$ws.Range("J55").Value = 58688
$ws.Range("L55").Value = 58688
$ws.Range("N55").Value = -59234
$ws.Range("H99").Value = 3791285.5
$ws.Range("I99").Value = 2267.6365
$ws.Range("J99").Value = 6997378
$ws.Range("K99").Value = 2267.6365
$ws.Range("L99").Value = 6997378
$ws.Range("M99").Value = -769.6365000000001
$ws.Range("N99").Value = -7000374
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8699.511
$ws.Range("I31").Value = 3489
$ws.Range("J31").Value = 11141.9375
$ws.Range("K31").Value = 3489
$ws.Range("L31").Value = 11141.9375
$ws.Range("M31").Value = -3194
$ws.Range("N31").Value = -11731.9375
$ws.Range("H34").Value = 8699.511
$ws.Range("I34").Value = 3489
$ws.Range("J34").Value = 11141.9375
$ws.Range("K34").Value = 3489
$ws.Range("L34").Value = 11141.9375
$ws.Range("M34").Value = -3287
$ws.Range("N34").Value = -11545.9375
$ws.Range("H56").Value = 65495
$ws.Range("J56").Value = 65495
$ws.Range("L56").Value = 65495
$ws.Range("N56").Value = -67185
$ws.Range("H59").Value = 99998.5
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H99").Value = 8624.3125
$ws.Range("I99").Value = 8887.666999999999
$ws.Range("J99").Value = 8285.714
$ws.Range("K99").Value = 8887.666999999999
$ws.Range("L99").Value = 8285.714
$ws.Range("M99").Value = -7389.666999999999
$ws.Range("N99").Value = -11281.714
$ws.Range("H107").Value = 1805.4783
$ws.Range("I107").Value = 1512.8572
$ws.Range("K107").Value = 1512.8572
$ws.Range("M107").Value = 407.1428000000001
$ws.Range("H126").Value = 8624.3125
$ws.Range("I126").Value = 8887.666999999999
$ws.Range("J126").Value = 8285.714
$ws.Range("K126").Value = 26663.001
$ws.Range("L126").Value = 24857.142
$ws.Range("M126").Value = -24193.001
$ws.Range("N126").Value = -29797.142
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 40009280
$ws.Range("J9").Value = 500
$ws.Range("L9").Value = 1500
$ws.Range("N9").Value = -1948
$ws.Range("H12").Value = 2381656.8
$ws.Range("J12").Value = 3333995.2
$ws.Range("L12").Value = 10001985.6
$ws.Range("N12").Value = -10002331.6
$ws.Range("H80").Value = 37041260
$ws.Range("J80").Value = 71433570
$ws.Range("L80").Value = 214300710
$ws.Range("N80").Value = -214302582
$ws.Range("H83").Value = 37041260
$ws.Range("J83").Value = 71433570
$ws.Range("L83").Value = 642902130
$ws.Range("N83").Value = -642911490
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("H132").Value = 11749.1
$ws.Range("I132").Value = 4798.8
$ws.Range("K132").Value = 43189.2
$ws.Range("M132").Value = -40659.2
$ws.Range("H141").Value = 8882.083000000001
$ws.Range("I141").Value = 3323.125
$ws.Range("K141").Value = 9969.375
$ws.Range("M141").Value = -4789.375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2320.5386
$ws.Range("I22").Value = 1075.5
$ws.Range("J22").Value = 2873.889
$ws.Range("K22").Value = 1075.5
$ws.Range("L22").Value = 2873.889
$ws.Range("M22").Value = -780.5
$ws.Range("N22").Value = -3463.889
$ws.Range("H27").Value = 2320.5386
$ws.Range("I27").Value = 1075.5
$ws.Range("J27").Value = 2873.889
$ws.Range("K27").Value = 1075.5
$ws.Range("L27").Value = 2873.889
$ws.Range("M27").Value = -968.5
$ws.Range("N27").Value = -3087.889
$ws.Range("H93").Value = 7325.4165
$ws.Range("I93").Value = 6790.5
$ws.Range("J93").Value = 10000
$ws.Range("K93").Value = 6790.5
$ws.Range("L93").Value = 10000
$ws.Range("M93").Value = -5542.5
$ws.Range("N93").Value = -12496
$ws.Range("H132").Value = 16676918
$ws.Range("I132").Value = 38469660
$ws.Range("K132").Value = 115408980
$ws.Range("M132").Value = -115406450
$ws.Range("H133").Value = 75333.336
$ws.Range("J133").Value = 75333.336
$ws.Range("L133").Value = 75333.336
$ws.Range("N133").Value = -80393.336
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 41667896
$ws.Range("I107").Value = 1337
$ws.Range("K107").Value = 4011
$ws.Range("M107").Value = -2091
$ws.Range("H132").Value = 14712247
$ws.Range("I132").Value = 27785570
$ws.Range("J132").Value = 4759.3125
$ws.Range("K132").Value = 83356710
$ws.Range("L132").Value = 14277.9375
$ws.Range("M132").Value = -83354180
$ws.Range("N132").Value = -19337.9375
